# Updates cryptocurrency "Price" (D) and "Volume(1h)" (E) columns on Sheet1
# to the latest scraped values. Values are textual (prices use "." as a
# thousands separator in this sheet, e.g. "55.508.39", and percentages keep
# their original padding, e.g. "  -3.66%  "), so each cell is forced to
# Text format before the write and restored to the default "Normal" style
# afterward so no formatting is left behind - only the displayed value
# changes, exactly like the upstream GitHub Actions scraper commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '55.508.39'
    'E2' = '  -3.66%  '
    'D3' = '2.912.86'
    'E3' = '  -3.88%  '
    'E4' = '  -0.10%  '
    'D5' = '500.30'
    'E5' = '  -2.11%  '
    'D6' = '133.24'
    'E6' = '  -4.47%  '
    'D7' = '1.00'
    'E7' = '  -0.05%  '
    'E8' = '  -4.84%  '
    'D9' = '7.16'
    'E9' = '  -4.18%  '
    'E10' = '  -5.53%  '
    'D11' = '0.350'
    'E11' = '  -4.47%  '
    'D12' = '3.403.75'
    'E12' = '  -4.27%  '
    'E13' = '  -3.82%  '
    'D14' = '25.70'
    'E14' = '  -3.33%  '
    'D15' = '0.0000158'
    'E15' = '  -4.31%  '
    'D16' = '55.468.04'
    'E16' = '  -3.84%  '
    'E17' = '  -2.98%  '
    'D18' = '2.909.05'
    'E18' = '  -4.19%  '
    'D19' = '12.53'
    'E19' = '  -2.79%  '
    'D20' = '7.66'
    'E20' = '  -4.05%  '
    'D21' = '314.74'
    'E21' = '  -5.60%  '
    'E22' = '  +0.16%  '
    'D23' = '0.483'
    'E23' = '  -2.98%  '
    'D24' = '62.30'
    'E24' = '  -3.41%  '
    'D25' = '3.023.91'
    'E25' = '  -4.37%  '
    'D26' = '0.999'
    'D27' = '0.160'
    'E27' = '  -4.47%  '
    'D28' = '0.0₃0848'
    'E28' = '  -8.12%  '
    'D29' = '6.35'
    'E29' = '  -6.41%  '
    'D30' = '6.90'
    'E30' = '  -6.89%  '
    'D31' = '1.76'
    'E31' = '  -2.78%  '
    'D32' = '19.69'
    'E32' = '  -5.21%  '
    'E33' = '  -7.21%  '
    'D34' = '149.24'
    'E34' = '  -4.40%  '
    'D35' = '4.37'
    'E35' = '  -6.95%  '
    'D36' = '5.59'
    'E36' = '  -4.41%  '
    'D37' = '24.59'
    'E37' = '  -0.18%  '
    'E38' = '  -7.25%  '
    'D39' = '0.0649'
    'E39' = '  -5.09%  '
    'E40' = '  -0.20%  '
    'D41' = '36.18'
    'E41' = '  -3.44%  '
    'D42' = '3.69'
    'E42' = '  -4.48%  '
    'D43' = '0.634'
    'E43' = '  -3.60%  '
    'D44' = '2.089.72'
    'E44' = '  -9.16%  '
    'E45' = '  -6.85%  '
    'D46' = '5.88'
    'E46' = '  -2.25%  '
    'D47' = '0.913'
    'E47' = '  -7.57%  '
    'E48' = '  -3.67%  '
    'D49' = '18.55'
    'E49' = '  -4.62%  '
    'D50' = '0.0834'
    'E50' = '  -6.61%  '
    'D51' = '1.68'
    'E51' = '  -8.08%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$cellRef]
    $cell.Style = 'Normal'
}
